# Insert four new "Title and Content" slides right after slide 1
# (the cover slide), pushing the existing slides down. This reproduces
# the target deck where slides 2-5 are brand new ("OOP?", "Pure OOP",
# "Hybrid OOP", "Why does Java does not support multiple inheritance?")
# and the previously-existing slides (Class Person, Class Pekerja,
# Class Caleg, Main Class part 1/2/3, RESULT) simply move down by four
# positions, unchanged.

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------
# Slide 2: "OOP?"
# ---------------------------------------------------------------
$s2 = $p.Slides.Add(2, 2)
$s2.Shapes.Item(1).TextFrame.TextRange.Text = "OOP?"

$body2 = $s2.Shapes.Item(2).TextFrame.TextRange
$body2.Text = "> OOP merupakan konsep pemrograman berbasis objek dimana orientasi nya berdasarkan apa yang ada dikehidupan sehari-hari"
$body2.ParagraphFormat.Bullet.Visible = 0
$rl2 = $s2.Shapes.Item(2).TextFrame.Ruler.Levels.Item(1)
$rl2.FirstMargin = 0
$rl2.LeftMargin = 0

# ---------------------------------------------------------------
# Slide 3: "Pure OOP"
# ---------------------------------------------------------------
$s3 = $p.Slides.Add(3, 2)
$s3.Shapes.Item(1).TextFrame.TextRange.Text = "Pure OOP"
$s3.Shapes.Item(2).TextFrame.TextRange.Text = "adalah sebuah bahasa yang mengharuskan program ditulis hanya berupa object saja. Contoh " + [char]0x2013 + " Eifel, Smaltalk, Ruby, Jade dan lain-lain."

# ---------------------------------------------------------------
# Slide 4: "Hybrid OOP"
# ---------------------------------------------------------------
$s4 = $p.Slides.Add(4, 2)
$s4.Shapes.Item(1).TextFrame.TextRange.Text = "Hybrid OOP"
$s4.Shapes.Item(2).TextFrame.TextRange.Text = "adalah bahasa yang dirancang untuk pemrograman object dengan beberapa elemen procedural."

# ---------------------------------------------------------------
# Slide 5: "Why does Java does not support multiple inheritance?"
# ---------------------------------------------------------------
$s5 = $p.Slides.Add(5, 2)
$s5.Shapes.Item(1).TextFrame.TextRange.Text = "Why does Java does not support multiple inheritance?"

$para1 = "Multiple infheritance adalah sebua class yang bisa menginherit/mewariskan atribut atau sifat-sifat dari dua atau lebih class."
$para2 = "Java, C#, etc tidak mensupport multiple inheritance karena dapat membuat ambiguitas. "
$para3 = "Alternatif dari multiple inheritance dengan menggunakan interface."
$s5.Shapes.Item(2).TextFrame.TextRange.Text = $para1 + "`r" + $para2 + "`r" + $para3

Write-Output ("Final slide count: " + $p.Slides.Count)
